$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table "表1" currently spans A1:Q16; add a new table row which will
# grow the table/autofilter ref (and worksheet dimension) to A1:Q17.
$tbl = $ws.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add()

# Fill in the new record (row 17) - only columns A-E are populated,
# matching the sparse pattern used by the other rows in this sheet.
# Script-path columns (D/E) are set before the Name column (B) so that
# the new shared-string entries are appended in the same order as the
# authoritative edit: "magnetic" first, then "地磁反转".
$ws.Range("A17").Value = 42000014
$ws.Range("D17").Value = "magnetic"
$ws.Range("E17").Value = "magnetic"
$ws.Range("B17").Value = "地磁反转"
$ws.Range("C17").Value = 0

# Match the author's final cursor position on the new row's Name cell.
$ws.Range("B17").Select() | Out-Null
